$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = "Aquari"
    6   = "Nexoon"
    28  = "Florian"
    36  = "Faerix"
    48  = "Zorblax"
    69  = "Faerix"
    94  = "Cybex"
    98  = "Faerix"
    112 = "Sentire"
    144 = "Florian"
    167 = "Emotivor"
    173 = "Emotivor"
    204 = "Emotivor"
    206 = "Mythron"
    214 = "Zorblax"
    223 = "Sentire"
    224 = "Nexoon"
    236 = "Faerix"
    241 = "Zorblax"
    271 = "Quixnar"
    277 = "Zorblax"
    282 = "Emotivor"
    292 = "Faerix"
    296 = "Mythron"
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
